$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "10/03/2025"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = "GPT5"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 231
$ws.Range("R3").Value = 258
$ws.Range("S3").Value = 0.9626512968299712
$ws.Range("T3").Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_2ndOctober_FewShotTest_Embeddings/ner_evaluation_results_GPT5_5_shot.txt"
$ws.Range("U3").Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_2ndOctober_FewShotTest_Embeddings/Stats/ner_evaluation_stats_GPT5_5_shot.txt"
$ws.Range("V3").Value = "4 MLGPU"
$ws.Range("W3").Value = "0.002 kWh"
